$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells carry textual (inline-string) values in the source data.
# Force text number-format first so Excel does not auto-coerce numeric-looking
# strings (e.g. "1.000", "102.60", "29.953.84") into actual numbers.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.953.84'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.61%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.876.12'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.73%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.91'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.78%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.02%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4922'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.71%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2911'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.80%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06622'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.71%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.877.90'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.65%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.70'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.21%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07240'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.67%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6643'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.91%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '86.16'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.77%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.893'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.35%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.937.80'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.72%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007840'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -6.11%  '

# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.06%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.74'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.72%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.121.19'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.63%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9997'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.06%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.765'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.00%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.744'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.63%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.044'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.98%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.04'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.20%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.83'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.28%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.98'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.06%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.913'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.24%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.394'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.01%  '

# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.01%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08721'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.09%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.953'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.36%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05046'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.27%  '

# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7095'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.86%  '

# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.111'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.16%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.670'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.67%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01781'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.45%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.685'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.54%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.165'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.97%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9288'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.19%  '

# Row 41
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9985'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.04%  '

# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4232'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.40%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.753'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.74%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.60'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.42%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.414'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.66%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1265'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.82%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05664'

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '32.46'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.87%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3771'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.81%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.253'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.75%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.90'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.45%  '
